# Update "想去人数" (F column) values on the 展览 (Exhibitions) sheet
# and the 全部类型 (All Types) sheet to reflect newly generated output.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# 展览 sheet (sheet1) — rows 4-29
$ws1.Range("F4").Value  = 74
$ws1.Range("F5").Value  = 1663
$ws1.Range("F6").Value  = 3264
$ws1.Range("F7").Value  = 826
$ws1.Range("F8").Value  = 2060
$ws1.Range("F9").Value  = 1972
$ws1.Range("F10").Value = 1018
$ws1.Range("F13").Value = 1611
$ws1.Range("F14").Value = 348
$ws1.Range("F18").Value = 71
$ws1.Range("F19").Value = 1443
$ws1.Range("F21").Value = 634
$ws1.Range("F23").Value = 10745
$ws1.Range("F24").Value = 11709
$ws1.Range("F29").Value = 443

# 全部类型 sheet (sheet4) — rows 6-35
$ws4.Range("F6").Value  = 74
$ws4.Range("F7").Value  = 1663
$ws4.Range("F8").Value  = 3264
$ws4.Range("F9").Value  = 826
$ws4.Range("F10").Value = 2060
$ws4.Range("F11").Value = 1972
$ws4.Range("F12").Value = 1018
$ws4.Range("F15").Value = 1611
$ws4.Range("F16").Value = 348
$ws4.Range("F22").Value = 71
$ws4.Range("F23").Value = 1443
$ws4.Range("F25").Value = 634
$ws4.Range("F27").Value = 10745
$ws4.Range("F28").Value = 11710
$ws4.Range("F35").Value = 443
